$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell updates that are plain text (not parseable as numbers) ---
$ws.Range('D2').Value = '41.379.05'
$ws.Range('E2').Value = '  -3.09%  '
$ws.Range('D3').Value = '2.469.75'
$ws.Range('E3').Value = '  -2.15%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '312.42'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.55%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '94.70'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -6.42%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.553'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -3.06%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.505'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -4.27%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '33.79'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -6.09%  '
$ws.Range('E11').Value = '  -2.64%  '
$ws.Range('E12').Value = '  -0.27%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '7.04'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -3.59%  '
$ws.Range('D14').Value = '2.846.04'
$ws.Range('E14').Value = '  -2.26%  '
$ws.Range('D15').Value = '2.473.49'
$ws.Range('E15').Value = '  -1.31%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '14.64'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -6.57%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.790'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -2.30%  '
$ws.Range('D18').Value = '41.328.59'
$ws.Range('E18').Value = '  -3.14%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.35'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -5.56%  '
$ws.Range('D20').Value = '0.0₃0922'
$ws.Range('E20').Value = '  -3.08%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.54'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -5.46%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '68.18'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.80%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '237.18'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -3.07%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.79'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -3.73%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.94'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -5.18%  '
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '24.53'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -6.53%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.24'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -4.09%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.74'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -4.11%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '36.25'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -7.37%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '153.62'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.61%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.61'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.98%  '
$ws.Range('E33').Value = '  -6.17%  '
$ws.Range('E34').Value = '  -1.08%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0758'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -3.72%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.03'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -5.36%  '
$ws.Range('E37').Value = '  -6.55%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '17.11'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -5.91%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '4.36'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.92%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.104'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -7.05%  '
$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.114'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -3.69%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '21.36'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.59%  '
$ws.Range('E43').Value = '  +0.24%  '
$ws.Range('D44').Value = '1.974.97'
$ws.Range('E44').Value = '  -0.66%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0287'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -4.38%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.09'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -6.23%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '8.69'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.70%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '70.18'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -3.21%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '76.26'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -4.98%  '
$ws.Range('E50').Value = '  -3.65%  '
$ws.Range('E51').Value = '  -5.37%  '
